$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number but must remain
# stored as text (matching the original inlineStr cell type in the sheet).
# We temporarily force a text number format so Excel does not auto-convert
# the assigned string into a numeric value, then clear the format again so
# no stray style index is left behind on the cell.
$textForceCells = @(
    'D5',
    'D7',
    'D10',
    'D12',
    'D14',
    'D15',
    'D16',
    'D19',
    'D20',
    'D22',
    'D25',
    'D26',
    'D27',
    'D29',
    'D32',
    'D33',
    'D35',
    'D36',
    'D37',
    'D41',
    'D42',
    'D44',
    'D45',
)
foreach ($c in $textForceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '38.151.04'
$ws.Range('E2').Value = '  +2.86%  '
$ws.Range('D3').Value = '2.059.04'
$ws.Range('E3').Value = '  +2.36%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '230.18'
$ws.Range('E5').Value = '  +1.81%  '
$ws.Range('E6').Value = '  +2.84%  '
$ws.Range('D7').Value = '59.53'
$ws.Range('E7').Value = '  +8.56%  '
$ws.Range('E9').Value = '  +3.28%  '
$ws.Range('D10').Value = '0.0814'
$ws.Range('E10').Value = '  +4.24%  '
$ws.Range('E11').Value = '  +2.14%  '
$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').Value = '14.78'
$ws.Range('E12').Value = '  +5.37%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '2.362.45'
$ws.Range('E13').Value = '  +2.28%  '
$ws.Range('D14').Value = '21.13'
$ws.Range('E14').Value = '  +6.80%  '
$ws.Range('D15').Value = '0.755'
$ws.Range('E15').Value = '  +2.56%  '
$ws.Range('D16').Value = '5.30'
$ws.Range('E16').Value = '  +1.91%  '
$ws.Range('D17').Value = '2.057.65'
$ws.Range('E17').Value = '  +2.41%  '
$ws.Range('D18').Value = '38.008.85'
$ws.Range('E18').Value = '  +2.75%  '
$ws.Range('D19').Value = '6.32'
$ws.Range('E19').Value = '  +1.55%  '
$ws.Range('D20').Value = '69.92'
$ws.Range('E20').Value = '  +2.50%  '
$ws.Range('D21').Value = '0.0₃0838'
$ws.Range('E21').Value = '  +3.02%  '
$ws.Range('D22').Value = '224.47'
$ws.Range('E22').Value = '  +0.80%  '
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').Value = '2.26'
$ws.Range('E25').Value = '  +4.07%  '
$ws.Range('D26').Value = '9.31'
$ws.Range('E26').Value = '  +3.66%  '
$ws.Range('D27').Value = '166.56'
$ws.Range('E27').Value = '  +1.23%  '
$ws.Range('E28').Value = '  +7.17%  '
$ws.Range('D29').Value = '19.06'
$ws.Range('E29').Value = '  +2.69%  '
$ws.Range('E30').Value = '  +2.83%  '
$ws.Range('E31').Value = '  +2.43%  '
$ws.Range('D32').Value = '4.57'
$ws.Range('E32').Value = '  +2.48%  '
$ws.Range('D33').Value = '4.61'
$ws.Range('E33').Value = '  +2.74%  '
$ws.Range('E34').Value = '  +11.04%  '
$ws.Range('D35').Value = '0.0608'
$ws.Range('E35').Value = '  +1.34%  '
$ws.Range('D36').Value = '2.33'
$ws.Range('E36').Value = '  +0.65%  '
$ws.Range('D37').Value = '6.10'
$ws.Range('E37').Value = '  +14.44%  '
$ws.Range('E38').Value = '  +5.29%  '
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('D40').Value = '1.536.43'
$ws.Range('E40').Value = '  +5.57%  '
$ws.Range('D41').Value = '98.38'
$ws.Range('E41').Value = '  +3.85%  '
$ws.Range('D42').Value = '0.0218'
$ws.Range('E42').Value = '  +2.61%  '
$ws.Range('E43').Value = '  +4.60%  '
$ws.Range('D44').Value = '16.86'
$ws.Range('E44').Value = '  +6.05%  '
$ws.Range('D45').Value = '0.0925'
$ws.Range('E45').Value = '  +2.20%  '
$ws.Range('E46').Value = '  +1.22%  '
$ws.Range('E47').Value = '  +13.49%  '
$ws.Range('E48').Value = '  +2.65%  '
$ws.Range('E49').Value = '  +2.71%  '
$ws.Range('E50').Value = '  +0.25%  '
$ws.Range('D51').Value = '2.249.97'
$ws.Range('E51').Value = '  +2.39%  '

foreach ($c in $textForceCells) {
    $ws.Range($c).ClearFormats()
}
